$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.31978701593237
$ws.Range("C2").Value = 10.37321111750963
$ws.Range("D2").Value = 7.073850283782638
$ws.Range("F2").Value = 35.44839233106196
$ws.Range("G2").Value = 48.24358697157206
$ws.Range("H2").Value = 19.00721503665709
$ws.Range("K2").Value = 9.347028763936462
$ws.Range("L2").Value = 11.03307850501169
$ws.Range("M2").Value = 15.27909655777725
$ws.Range("B3").Value = 13.16507134068635
$ws.Range("C3").Value = 10.35896940173343
$ws.Range("D3").Value = 7.069099978129377
$ws.Range("F3").Value = 35.34764491734436
$ws.Range("G3").Value = 48.04538322495618
$ws.Range("H3").Value = 19.015798473328
$ws.Range("K3").Value = 9.239720896800051
$ws.Range("L3").Value = 11.04026037074701
$ws.Range("M3").Value = 15.26785055887856
$ws.Range("B4").Value = 13.07289434882427
$ws.Range("C4").Value = 10.34985501354155
$ws.Range("D4").Value = 7.06594675242521
$ws.Range("F4").Value = 35.29274539712124
$ws.Range("G4").Value = 47.93353176218789
$ws.Range("H4").Value = 19.02379531272618
$ws.Range("K4").Value = 9.175852918261226
$ws.Range("L4").Value = 11.04629671396224
$ws.Range("M4").Value = 15.26374542544712
$ws.Range("B5").Value = 13.03608688490337
$ws.Range("C5").Value = 10.34604770174722
$ws.Range("D5").Value = 7.064602318643523
$ws.Range("F5").Value = 35.27213575718156
$ws.Range("G5").Value = 47.89045450618994
$ws.Range("H5").Value = 19.02773964094423
$ws.Range("K5").Value = 9.150366656093398
$ws.Range("L5").Value = 11.04916596950729
$ws.Range("M5").Value = 15.26277915614193
$ws.Range("B6").Value = 13.03002198796123
$ws.Range("C6").Value = 10.34540987295567
$ws.Range("D6").Value = 7.064375484873548
$ws.Range("F6").Value = 35.26882030720304
$ws.Range("G6").Value = 47.88345347231627
$ws.Range("H6").Value = 19.02843599361254
$ws.Range("K6").Value = 9.146168256477438
$ws.Range("L6").Value = 11.04966714244239
$ws.Range("M6").Value = 15.26266145154418
$ws.Range("B7").Value = 13.07239483153562
$ws.Range("C7").Value = 10.34980404399695
$ws.Range("D7").Value = 7.065928861711559
$ws.Range("F7").Value = 35.29246029779563
$ws.Range("G7").Value = 47.93294063885499
$ws.Range("H7").Value = 19.02384573181144
$ws.Range("K7").Value = 9.175506971007675
$ws.Range("L7").Value = 11.0463337517074
$ws.Range("M7").Value = 15.26372952992704
$ws.Range("B8").Value = 13.26588432090421
$ws.Range("C8").Value = 10.36837728392232
$ws.Range("D8").Value = 7.072261470590798
$ws.Range("F8").Value = 35.4122180031178
$ws.Range("G8").Value = 48.173219952533
$ws.Range("H8").Value = 19.00960871331594
$ws.Range("K8").Value = 9.309629767159432
$ws.Range("L8").Value = 11.03521744536577
$ws.Range("M8").Value = 15.27463943886822
$ws.Range("B9").Value = 13.66557428014141
$ws.Range("C9").Value = 10.40186984811006
$ws.Range("D9").Value = 7.08280287976408
$ws.Range("F9").Value = 35.70166996508924
$ws.Range("G9").Value = 48.72119724626803
$ws.Range("H9").Value = 19.00332215487769
$ws.Range("K9").Value = 9.587187032531007
$ws.Range("L9").Value = 11.02630646812299
$ws.Range("M9").Value = 15.31812374386322
$ws.Range("B10").Value = 13.96868599558538
$ws.Range("C10").Value = 10.42469993061927
$ws.Range("D10").Value = 7.089403453867901
$ws.Range("F10").Value = 35.9466871539723
$ws.Range("G10").Value = 49.16861943110901
$ws.Range("H10").Value = 19.01188091924701
$ws.Range("K10").Value = 9.797945181483096
$ws.Range("L10").Value = 11.02758731961967
$ws.Range("M10").Value = 15.3633468812244
$ws.Range("B11").Value = 14.10801193545854
$ws.Range("C11").Value = 10.43470210063787
$ws.Range("D11").Value = 7.09215779786712
$ws.Range("F11").Value = 36.06494351241076
$ws.Range("G11").Value = 49.38141169154409
$ws.Range("H11").Value = 19.01862983899846
$ws.Range("K11").Value = 9.894872565591735
$ws.Range("L11").Value = 11.02986144049537
$ws.Range("M11").Value = 15.38675373591897
$ws.Range("B12").Value = 14.16092627282678
$ws.Range("C12").Value = 10.43843466860973
$ws.Range("D12").Value = 7.093165092490661
$ws.Range("F12").Value = 36.11067901375388
$ws.Range("G12").Value = 49.46327666852839
$ws.Range("H12").Value = 19.02159515531886
$ws.Range("K12").Value = 9.931691471607893
$ws.Range("L12").Value = 11.03096490123005
$ws.Range("M12").Value = 15.39602009452258
$ws.Range("B13").Value = 14.14952417117313
$ws.Range("C13").Value = 10.43763323966404
$ws.Range("D13").Value = 7.092949742200183
$ws.Range("F13").Value = 36.10078699097348
$ws.Range("G13").Value = 49.4455891819842
$ws.Range("H13").Value = 19.02093831997463
$ws.Range("K13").Value = 9.923757343678009
$ws.Range("L13").Value = 11.0307164908435
$ws.Range("M13").Value = 15.39400659119252
$ws.Range("B14").Value = 14.11236242780376
$ws.Range("C14").Value = 10.43501028375731
$ws.Range("D14").Value = 7.0922413844593
$ws.Range("F14").Value = 36.06868720945619
$ws.Range("G14").Value = 49.38812128829041
$ws.Range("H14").Value = 19.01886559374954
$ws.Range("K14").Value = 9.897899590850326
$ws.Range("L14").Value = 11.02994737253717
$ws.Range("M14").Value = 15.38750804389368
$ws.Range("B15").Value = 14.08961838431569
$ws.Range("C15").Value = 10.43339648277012
$ws.Range("D15").Value = 7.091802843333207
$ws.Range("F15").Value = 36.04914874702753
$ws.Range("G15").Value = 49.35308650194785
$ws.Range("H15").Value = 19.01764930472201
$ws.Range("K15").Value = 9.882074814232091
$ws.Range("L15").Value = 11.02950779053099
$ws.Range("M15").Value = 15.38357978541365
$ws.Range("B16").Value = 13.95960518768604
$ws.Range("C16").Value = 10.42403855884475
$ws.Range("D16").Value = 7.089218453720339
$ws.Range("F16").Value = 35.93909366232432
$ws.Range("G16").Value = 49.15489536164658
$ws.Range("H16").Value = 19.01149724336276
$ws.Range("K16").Value = 9.791628776700504
$ws.Range("L16").Value = 11.02747266010136
$ws.Range("M16").Value = 15.36187381608932
$ws.Range("B17").Value = 13.8801774427313
$ws.Range("C17").Value = 10.41819965315328
$ws.Range("D17").Value = 7.087569412320128
$ws.Range("F17").Value = 35.87330339821581
$ws.Range("G17").Value = 49.0356506304466
$ws.Range("H17").Value = 19.00845389241189
$ws.Range("K17").Value = 9.736386349922256
$ws.Range("L17").Value = 11.02665677129639
$ws.Range("M17").Value = 15.34928054320089
$ws.Range("B18").Value = 13.83463123886274
$ws.Range("C18").Value = 10.41480518104753
$ws.Range("D18").Value = 7.086597557612794
$ws.Range("F18").Value = 35.83610380322587
$ws.Range("G18").Value = 48.96793851026416
$ws.Range("H18").Value = 19.00697228171916
$ws.Range("K18").Value = 9.704713605469278
$ws.Range("L18").Value = 11.02634676419727
$ws.Range("M18").Value = 15.34230435371124
$ws.Range("B19").Value = 13.81923547375309
$ws.Range("C19").Value = 10.41364967029715
$ws.Range("D19").Value = 7.086264490941209
$ws.Range("F19").Value = 35.82361946555094
$ws.Range("G19").Value = 48.94516387449666
$ws.Range("H19").Value = 19.0065168341106
$ws.Range("K19").Value = 9.694008272538976
$ws.Range("L19").Value = 11.02626918241603
$ws.Range("M19").Value = 15.33998835891455
$ws.Range("B20").Value = 13.88861870486304
$ws.Range("C20").Value = 10.41882494897187
$ws.Range("D20").Value = 7.087747374349332
$ws.Range("F20").Value = 35.88024068133108
$ws.Range("G20").Value = 49.04825427529602
$ws.Range("H20").Value = 19.00875004424494
$ws.Range("K20").Value = 9.742256790301223
$ws.Range("L20").Value = 11.02672714651167
$ws.Range("M20").Value = 15.35059350627482
$ws.Range("B21").Value = 14.1232739492237
$ws.Range("C21").Value = 10.43578220245908
$ws.Range("D21").Value = 7.092450415696105
$ws.Range("F21").Value = 36.07808998245768
$ws.Range("G21").Value = 49.40496650178197
$ws.Range("H21").Value = 19.01946329478696
$ws.Range("K21").Value = 9.905491819481902
$ws.Range("L21").Value = 11.03016671322972
$ws.Range("M21").Value = 15.38940593703677
$ws.Range("B22").Value = 14.27751232332937
$ws.Range("C22").Value = 10.44654399372238
$ws.Range("D22").Value = 7.095315942611971
$ws.Range("F22").Value = 36.21294709976146
$ws.Range("G22").Value = 49.64556647369531
$ws.Range("H22").Value = 19.02885215592837
$ws.Range("K22").Value = 10.0128267332902
$ws.Range("L22").Value = 11.03382647003686
$ws.Range("M22").Value = 15.41711673566057
$ws.Range("B23").Value = 14.19512917272679
$ws.Range("C23").Value = 10.44082954261784
$ws.Range("D23").Value = 7.093805609138085
$ws.Range("F23").Value = 36.14047137222486
$ws.Range("G23").Value = 49.51648625827615
$ws.Range("H23").Value = 19.02362309768704
$ws.Range("K23").Value = 9.955492465232201
$ws.Range("L23").Value = 11.0317443454543
$ws.Range("M23").Value = 15.4021141544556
$ws.Range("B24").Value = 13.88480203911873
$ws.Range("C24").Value = 10.41854236989145
$ws.Range("D24").Value = 7.087666991823347
$ws.Range("F24").Value = 35.87710238977251
$ws.Range("G24").Value = 49.04255353716059
$ws.Range("H24").Value = 19.00861531890383
$ws.Range("K24").Value = 9.739602490825748
$ws.Range("L24").Value = 11.02669483435821
$ws.Range("M24").Value = 15.34999909342383
$ws.Range("B25").Value = 13.55557520074334
$ws.Range("C25").Value = 10.39312192021396
$ws.Range("D25").Value = 7.080152087152469
$ws.Range("F25").Value = 35.61760884469169
$ws.Range("G25").Value = 48.56492744620828
$ws.Range("H25").Value = 19.00270719330184
$ws.Range("K25").Value = 9.510751537907359
$ws.Range("L25").Value = 11.02733995020547
$ws.Range("M25").Value = 15.30401336444766
